$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the "student type" row entirely (originally row 10)
$ws.Rows.Item(10).Delete()

# Delete the "spend time per week" row (originally row 12, now row 11 after the first delete)
$ws.Rows.Item(11).Delete()

# Set column F width to match target
$ws.Columns.Item(6).ColumnWidth = 34.6640625

# Update sheet view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("D3").Select()
